# Insert one new data row above current row 51 (pushes existing rows 51-122
# down to 52-123) and populate it with the new daily price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51:51").Insert()

$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44579
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112045
$ws.Range("G51").Value = "Zapallo"
$ws.Range("H51").Value = "Camote"
$ws.Range("I51").Value = "1a nueva(o)"
$ws.Range("J51").Value = 200
$ws.Range("K51").Value = 300
$ws.Range("L51").Value = 350
$ws.Range("M51").Value = 325
$ws.Range("N51").Value = "$/kilo (volumen en unidades)"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 325
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"

"done"
